# Update the "Ensemble" row (row 6) metrics to reflect the improved
# ensemble model that weights individual classifiers by their accuracy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B6" = 0.8261496225120111
    "C6" = 0.838139974546402
    "D6" = 0.8261496225120111
    "E6" = 0.8265545525823867
    "F6" = 0.8518874399450928
    "G6" = 0.8579622526349604
    "H6" = 0.8518874399450928
    "I6" = 0.8518664635506795
    "J6" = 0.8175245939144361
    "K6" = 0.8279543546961212
    "L6" = 0.8175245939144361
    "M6" = 0.8180653797168308
    "N6" = 0.8411576298329901
    "O6" = 0.8506178001618393
    "P6" = 0.8411576298329901
    "Q6" = 0.8407569444655187
    "R6" = 0.8497140242507436
    "S6" = 0.8545892027115745
    "T6" = 0.8497140242507436
    "U6" = 0.8495344564060954
    "V6" = 0.8561427590940289
    "W6" = 0.8624865978297163
    "X6" = 0.8561427590940289
    "Y6" = 0.8566615758221836
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
